$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to force numeric-looking strings to be stored as literal text
# (preserves trailing zeros / exact formatting instead of Excel coercing to a number).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = "42.728.87"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.360.92"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.04%  "
$helper.Value = "317.07"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -3.15%  "
$helper.Value = "109.27"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  -0.03%  "
$helper.Value = "0.620"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -4.46%  "
$helper.Value = "41.91"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  -1.24%  "
$helper.Value = "8.61"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.77%  "
$helper.Value = "1.01"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("E14").Value = "  +0.02%  "
$helper.Value = "16.11"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "2.714.29"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "2.451.03"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "42.697.64"
$ws.Range("E18").Value = "  -1.10%  "
$helper.Value = "7.72"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +0.06%  "
$helper.Value = "76.21"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  -2.35%  "
$helper.Value = "256.29"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -5.45%  "
$helper.Value = "2.33"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -3.95%  "
$helper.Value = "9.48"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -0.01%  "
$helper.Value = "11.47"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.65%  "
$helper.Value = "22.88"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  +2.19%  "
$helper.Value = "37.39"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.13%  "
$helper.Value = "171.91"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -2.24%  "
$helper.Value = "0.0893"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -3.88%  "
$helper.Value = "6.05"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +1.92%  "
$helper.Value = "2.90"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -8.63%  "
$helper.Value = "0.124"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +13.54%  "
$ws.Range("E36").Value = "  -2.50%  "
$helper.Value = "4.71"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -4.58%  "
$ws.Range("E38").Value = "  -0.44%  "
$helper.Value = "3.93"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -4.94%  "
$helper.Value = "2.70"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("E42").Value = "  -5.28%  "
$helper.Value = "71.17"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("E44").Value = "  +0.03%  "
$helper.Value = "12.22"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -1.16%  "
$helper.Value = "112.54"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -7.64%  "
$helper.Value = "5.54"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$helper.Value = "9.22"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$helper.Value = "85.93"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -3.64%  "
$helper.Value = "76.71"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("E51").Value = "  -1.25%  "

# Clean up helper cell and clipboard state
$helper.Clear()
$excel.CutCopyMode = $false
